$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text storage (matches source inline-string cells) so that
# numeric-looking values (e.g. "329.90", "21") are not auto-converted
# to Number/Percentage types by the .Value setter.
$ws.Range("D2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "329.90"
$ws.Range("E2").Value = "0.29%"
$ws.Range("G2").Value = "21"
$ws.Range("D3").Value = "41.16"
$ws.Range("E3").Value = "1.03%"
$ws.Range("G3").Value = "21"
$ws.Range("D4").Value = "5.695"
$ws.Range("E4").Value = "-0.71%"
$ws.Range("G4").Value = "21"
$ws.Range("D5").Value = "0.08070"
$ws.Range("E5").Value = "-0.48%"
$ws.Range("G5").Value = "21"
$ws.Range("D6").Value = "2.019"
$ws.Range("E6").Value = "2.83%"
$ws.Range("G6").Value = "21"
$ws.Range("D7").Value = "8.726"
$ws.Range("E7").Value = "-0.34%"
$ws.Range("G7").Value = "21"
$ws.Range("D8").Value = "4.524"
$ws.Range("E8").Value = "-1.39%"
$ws.Range("G8").Value = "21"
$ws.Range("G9").Value = "21"
$ws.Range("D10").Value = "0.9232"
$ws.Range("E10").Value = "-2.36%"
$ws.Range("G10").Value = "21"
$ws.Range("D11").Value = "0.1272"
$ws.Range("E11").Value = "-1.77%"
$ws.Range("G11").Value = "21"
$ws.Range("D12").Value = "0.1937"
$ws.Range("E12").Value = "-2.78%"
$ws.Range("G12").Value = "21"
$ws.Range("D13").Value = "8.254"
$ws.Range("E13").Value = "-7.83%"
$ws.Range("G13").Value = "21"
$ws.Range("D14").Value = "0.09379"
$ws.Range("E14").Value = "0.80%"
$ws.Range("G14").Value = "21"
$ws.Range("D15").Value = "0.03706"
$ws.Range("E15").Value = "5.27%"
$ws.Range("G15").Value = "21"
$ws.Range("D16").Value = "0.1054"
$ws.Range("G16").Value = "21"
$ws.Range("D17").Value = "0.001295"
$ws.Range("E17").Value = "-2.38%"
$ws.Range("G17").Value = "21"
$ws.Range("D18").Value = "0.006262"
$ws.Range("E18").Value = "2.42%"
$ws.Range("G18").Value = "21"
$ws.Range("D19").Value = "3.367"
$ws.Range("E19").Value = "-0.12%"
$ws.Range("G19").Value = "21"
$ws.Range("E20").Value = "-2.52%"
$ws.Range("G20").Value = "21"
$ws.Range("D21").Value = "0.1417"
$ws.Range("E21").Value = "-0.19%"
$ws.Range("G21").Value = "21"
$ws.Range("D22").Value = "0.2654"
$ws.Range("E22").Value = "10.14%"
$ws.Range("G22").Value = "21"
$ws.Range("D23").Value = "0.04419"
$ws.Range("E23").Value = "-0.47%"
$ws.Range("G23").Value = "21"
$ws.Range("D24").Value = "0.001259"
$ws.Range("E24").Value = "-0.05%"
$ws.Range("G24").Value = "21"
$ws.Range("D25").Value = "0.004342"
$ws.Range("E25").Value = "-0.63%"
$ws.Range("G25").Value = "21"
$ws.Range("D26").Value = "0.0001242"
$ws.Range("E26").Value = "13.53%"
$ws.Range("G26").Value = "21"
$ws.Range("G27").Value = "21"
$ws.Range("G28").Value = "21"
$ws.Range("G29").Value = "21"
$ws.Range("G30").Value = "21"
$ws.Range("G31").Value = "21"
$ws.Range("G32").Value = "21"
$ws.Range("G33").Value = "21"
$ws.Range("G34").Value = "21"
$ws.Range("G35").Value = "21"
$ws.Range("G36").Value = "21"
$ws.Range("G37").Value = "21"
$ws.Range("G38").Value = "21"
$ws.Range("D39").Value = "0.02864"
$ws.Range("E39").Value = "16.21%"
$ws.Range("G39").Value = "21"
$ws.Range("D40").Value = "0.05464"
$ws.Range("E40").Value = "2.76%"
$ws.Range("G40").Value = "21"
$ws.Range("D41").Value = "0.007716"
$ws.Range("E41").Value = "3.16%"
$ws.Range("G41").Value = "21"
$ws.Range("D42").Value = "0.009948"
$ws.Range("E42").Value = "12.38%"
$ws.Range("G42").Value = "21"
$ws.Range("D43").Value = "0.1418"
$ws.Range("E43").Value = "-1.38%"
$ws.Range("G43").Value = "21"
$ws.Range("D44").Value = "0.002134"
$ws.Range("E44").Value = "0.27%"
$ws.Range("G44").Value = "21"
$ws.Range("D45").Value = "0.01184"
$ws.Range("E45").Value = "13.49%"
$ws.Range("G45").Value = "21"
$ws.Range("D46").Value = "0.00006764"
$ws.Range("E46").Value = "-1.64%"
$ws.Range("G46").Value = "21"
$ws.Range("E47").Value = "0.04%"
$ws.Range("G47").Value = "21"
$ws.Range("D48").Value = "0.003002"
$ws.Range("E48").Value = "-14.67%"
$ws.Range("G48").Value = "21"
$ws.Range("D49").Value = "0.002282"
$ws.Range("E49").Value = "34.10%"
$ws.Range("G49").Value = "21"
$ws.Range("E50").Value = "0.04%"
$ws.Range("G50").Value = "21"
$ws.Range("E51").Value = "0.04%"
$ws.Range("G51").Value = "21"

# Restore default (General/no explicit number format) styling so the
# written cells keep Text type without leaving a style-index footprint.
$ws.Range("D2:G51").Style = "Normal"
